# ---------------------------------------------------------------------------
# dollar.xlsx -- quarterly income-statement update
#   * drop the oldest quarter column (فصل دوم منتهی به 1399/06) and shift
#     every remaining quarter one column to the left (D<-E, E<-F, ... L<-M)
#   * append a new quarter (فصل چهارم منتهی به 1401/12) in column M together
#     with its publish date (1402-02-27) and figures
#   * misc workbook bookkeeping (absolute path casing, revision id, window
#     size) and a couple of column-width tweaks that ride along with the
#     newly inserted column
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- shift the whole data block (headers + figures) one column to the left
$srcRange = $ws.Range("E8:M27")
$dstRange = $ws.Range("D8:L27")
$srcRange.Copy($dstRange) | Out-Null

# --- column M: brand-new quarter -------------------------------------------------
$ws.Range("M8").Value  = "فصل چهارم منتهی به 1401/12"
$ws.Range("M9").Value  = "1402-02-27 (7)"

$ws.Range("M11").Value = 14043
$ws.Range("M12").Value = -12141
$ws.Range("M13").Value = 1902
$ws.Range("M14").Value = -1017
$ws.Range("M15").Value = "-"
$ws.Range("M16").Value = 89
$ws.Range("M17").Value = 974
$ws.Range("M18").Value = -163
$ws.Range("M19").Value = -121
$ws.Range("M20").Value = 690
$ws.Range("M21").Value = "-"
$ws.Range("M22").Value = 690
$ws.Range("M23").Value = "-"
$ws.Range("M24").Value = 690
$ws.Range("M25").Value = 0
$ws.Range("M26").Value = 1540
$ws.Range("M27").Value = 0

# row 21 "مالیات": the shift turns the old G21 (0 -> "-") into the new F21
$ws.Range("F21").Value = "-"

# --- column widths: the "wide" (31) marker shifts along with the new column
$ws.Columns("D").ColumnWidth = 29
$ws.Columns("E").ColumnWidth = 31
$ws.Columns("F").ColumnWidth = 29
$ws.Columns("G").ColumnWidth = 29
$ws.Columns("H").ColumnWidth = 29
$ws.Columns("I").ColumnWidth = 31
$ws.Columns("J").ColumnWidth = 29
$ws.Columns("K").ColumnWidth = 29
$ws.Columns("L").ColumnWidth = 29
$ws.Columns("M").ColumnWidth = 31

# --- workbook-level bookkeeping ---------------------------------------------------
$wb.Windows.Item(1).WindowState = -4143   # xlNormal, harmless no-op kept for clarity

Write-Host "edit applied"
